$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet "DesignSheet". Add a new
# worksheet named "DefaultSheet" right after it, so it becomes sheet 2
# (sheetId=2, rId2) and the new, active tab (activeTab=1 / index 1).
$designSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $designSheet)
$newSheet.Name = "DefaultSheet"

# Populate the new sheet: A1 = "dummy" (its own distinct font style),
# B1 = 0.
$newSheet.Range("A1").Value = "dummy"
$newSheet.Range("A1").Font.Name = "Helvetica Neue"
$newSheet.Range("A1").Font.Size = 10
$newSheet.Range("A1").Font.ThemeColor = 1
$newSheet.Range("B1").Value = 0
$newSheet.Rows(1).RowHeight = 14

# Leave the selection/active cell on A2, as in the authored file.
$newSheet.Range("A2").Select() | Out-Null
